# Applies the symbol-list update described in the commit:
# "Updated symbol list on Fri Jan 27 19:46:17 UTC 2023 with GitHub Actions"
#
# Columns D (Price) and E (Volume(1h)) hold numeric-looking values that are
# stored as literal text in the workbook (e.g. "308.42", "1.19%"). A leading
# apostrophe forces Excel to keep them as text instead of coercing them to
# numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''308.42'
$ws.Range("E2").Value = '''1.19%'
$ws.Range("E3").Value = '''1.30%'
$ws.Range("D4").Value = '''5.060'
$ws.Range("E4").Value = '''0.72%'
$ws.Range("D5").Value = '''0.08121'
$ws.Range("E5").Value = '''0.45%'
$ws.Range("D6").Value = '''2.042'
$ws.Range("E6").Value = '''4.79%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = '''7.875'
$ws.Range("E7").Value = '''0.46%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9267'
$ws.Range("E8").Value = '''-0.55%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.1413'
$ws.Range("E9").Value = '''12.58%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1927'
$ws.Range("E10").Value = '''0.82%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.09115'
$ws.Range("E11").Value = '''-1.28%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03450'
$ws.Range("E12").Value = '''-1.61%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09924'
$ws.Range("E13").Value = '''-0.07%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001406'
$ws.Range("E14").Value = '''-0.68%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006216'
$ws.Range("E15").Value = '''-7.20%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.836'
$ws.Range("E16").Value = '''6.13%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''4.157'
$ws.Range("E17").Value = '''0.39%'
$ws.Range("D18").Value = '''3.487'
$ws.Range("E18").Value = '''12.99%'
$ws.Range("D19").Value = '''0.3441'
$ws.Range("E19").Value = '''-0.03%'
$ws.Range("D20").Value = '''0.1292'
$ws.Range("D21").Value = '''4.795'
$ws.Range("E21").Value = '''-7.40%'
$ws.Range("E22").Value = '''-7.51%'
$ws.Range("D23").Value = '''0.04388'
$ws.Range("E23").Value = '''-0.43%'
$ws.Range("D24").Value = '''0.001233'
$ws.Range("E24").Value = '''-0.12%'
$ws.Range("E25").Value = '''4.10%'
$ws.Range("E27").Value = '''-0.03%'
$ws.Range("D39").Value = '''0.02038'
$ws.Range("E39").Value = '''3.54%'
$ws.Range("D40").Value = '''0.05146'
$ws.Range("E40").Value = '''-0.41%'
$ws.Range("D41").Value = '''0.007488'
$ws.Range("E41").Value = '''-1.19%'
$ws.Range("D42").Value = '''0.01010'
$ws.Range("E42").Value = '''-0.36%'
$ws.Range("E43").Value = '''0.14%'
$ws.Range("E44").Value = '''1.39%'
$ws.Range("D45").Value = '''0.009733'
$ws.Range("E45").Value = '''-9.00%'
$ws.Range("D46").Value = '''0.00006290'
$ws.Range("E46").Value = '''-1.31%'
$ws.Range("E47").Value = '''-0.05%'
$ws.Range("E48").Value = '''2.03%'
$ws.Range("E49").Value = '''-22.00%'
$ws.Range("E50").Value = '''-0.05%'
$ws.Range("E51").Value = '''-0.05%'
